# Auto-generated edit script applying the diff to the cryptos worksheet
# Updates Coin (B), Link (C), Price (D), Volume(1h) (E) columns for the
# rows that changed per the latest GitHub Actions refresh of the cryptos list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.854.41"
$ws.Range("E2").Value = "  +2.67%  "

$ws.Range("D3").Value = "3.701.14"
$ws.Range("E3").Value = "  +5.50%  "

$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "'419.27"
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("D6").Value = "'131.12"
$ws.Range("E6").Value = "  -1.65%  "

$ws.Range("D7").Value = "3.693.92"
$ws.Range("E7").Value = "  +5.60%  "

$ws.Range("D8").Value = "'0.645"
$ws.Range("E8").Value = "  -0.52%  "

$ws.Range("D10").Value = "'0.775"
$ws.Range("E10").Value = "  -1.52%  "

$ws.Range("D11").Value = "'0.185"
$ws.Range("E11").Value = "  +14.48%  "

$ws.Range("D12").Value = "'0.0000406"
$ws.Range("E12").Value = "  +58.57%  "

$ws.Range("D13").Value = "'43.16"
$ws.Range("E13").Value = "  -1.49%  "

$ws.Range("D14").Value = "'10.58"
$ws.Range("E14").Value = "  +5.72%  "

$ws.Range("D15").Value = "4.286.00"
$ws.Range("E15").Value = "  +5.41%  "

$ws.Range("E16").Value = "  -0.83%  "

$ws.Range("D17").Value = "'20.79"
$ws.Range("E17").Value = "  +1.11%  "

$ws.Range("D18").Value = "3.718.28"
$ws.Range("E18").Value = "  +5.83%  "

$ws.Range("D19").Value = "'13.44"
$ws.Range("E19").Value = "  +8.65%  "

$ws.Range("D20").Value = "'1.14"
$ws.Range("E20").Value = "  +3.40%  "

$ws.Range("D21").Value = "66.824.93"
$ws.Range("E21").Value = "  +2.82%  "

$ws.Range("D22").Value = "'446.38"
$ws.Range("E22").Value = "  -3.13%  "

$ws.Range("D23").Value = "'16.27"
$ws.Range("E23").Value = "  +21.84%  "

$ws.Range("D24").Value = "'89.77"
$ws.Range("E24").Value = "  -0.97%  "

$ws.Range("E25").Value = "  -1.83%  "

$ws.Range("D26").Value = "'37.46"
$ws.Range("E26").Value = "  +10.33%  "

$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").Value = "'10.19"
$ws.Range("E27").Value = "  +2.03%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'3.32"
$ws.Range("E28").Value = "  -1.20%  "

$ws.Range("D29").Value = "'5.10"
$ws.Range("E29").Value = "  +5.69%  "

$ws.Range("D30").Value = "'0.125"
$ws.Range("E30").Value = "  +9.23%  "

$ws.Range("D31").Value = "'12.73"
$ws.Range("E31").Value = "  +2.00%  "

$ws.Range("D32").Value = "'2.71"
$ws.Range("E32").Value = "  +0.77%  "

$ws.Range("D33").Value = "'7.31"
$ws.Range("E33").Value = "  -3.59%  "

$ws.Range("D34").Value = "'0.164"
$ws.Range("E34").Value = "  +0.87%  "

$ws.Range("D35").Value = "'42.00"
$ws.Range("E35").Value = "  +3.92%  "

$ws.Range("D36").Value = "'57.38"
$ws.Range("E36").Value = "  -0.72%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").Value = "'0.0496"
$ws.Range("E38").Value = "  -1.21%  "

$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0745"
$ws.Range("E39").Value = "  +5.56%  "

$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").Value = "'3.10"
$ws.Range("E40").Value = "  +32.09%  "

$ws.Range("E41").Value = "  +2.61%  "

$ws.Range("D42").Value = "'29.39"
$ws.Range("E42").Value = "  +35.74%  "

$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "'3.43"
$ws.Range("E43").Value = "  +40.73%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("B45").Value = "LidoDAOToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D45").Value = "'3.43"
$ws.Range("E45").Value = "  +3.38%  "

$ws.Range("D46").Value = "'148.87"
$ws.Range("E46").Value = "  +1.63%  "

$ws.Range("E47").Value = "  +4.37%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'2.66"
$ws.Range("E48").Value = "  -3.77%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'4.38"
$ws.Range("E49").Value = "  -3.02%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.90"
$ws.Range("E50").Value = "  -7.47%  "

$ws.Range("D51").Value = "'0.309"
$ws.Range("E51").Value = "  -3.10%  "

Write-Output "Applied updates to cryptos worksheet"
